# tests/test_data/test_dfs_tax_sched.xlsx
# "tax sched test correction #3"
#
# 1) Remove the TRUNC() wrapper from the L/M/N/O tax-schedule formulas in
#    rows 2-6 so the cells return the raw (un-truncated) computed value.
# 2) Row 4's formula additionally picks up corrected bracket coefficients
#    (912.17 / 1038 / 8172 / 15694 instead of 939.68 / 1007 / 8064 / 15576),
#    and its R4/S4 threshold inputs move from 8004/13469 to 8005/13470.
# 3) Move the sheet's view/selection from G6 to N17 (with the viewport
#    scrolled so column J is the leftmost visible column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Formula = '=(C2>9000)*(C2<13996)*(997.8*(C2-9000)/10000+1400)*(C2-9000)/10000+(C2>13996)*(C2<54949)*((220.13*(C2-13996)/10000+2397)*(C2-13996)/10000+948.49)+(C2>54950)*(C2<260532)*(0.42*C2-8621.75)+(C2>260532)*(0.45*C2-16437.7)'
$ws.Range("M2").Formula = '=(D2>9000)*(D2<13996)*(997.8*(D2-9000)/10000+1400)*(D2-9000)/10000+(D2>13996)*(D2<54949)*((220.13*(D2-13996)/10000+2397)*(D2-13996)/10000+948.49)+(D2>54950)*(D2<260532)*(0.42*D2-8621.75)+(D2>260532)*(0.45*D2-16437.7)'
$ws.Range("N2").Formula = '=(E2>9000)*(E2<13996)*(997.8*(E2-9000)/10000+1400)*(E2-9000)/10000+(E2>13996)*(E2<54949)*((220.13*(E2-13996)/10000+2397)*(E2-13996)/10000+948.49)+(E2>54950)*(E2<260532)*(0.42*E2-8621.75)+(E2>260532)*(0.45*E2-16437.7)'
$ws.Range("O2").Formula = '=(F2>9000)*(F2<13996)*(997.8*(F2-9000)/10000+1400)*(F2-9000)/10000+(F2>13996)*(F2<54949)*((220.13*(F2-13996)/10000+2397)*(F2-13996)/10000+948.49)+(F2>54950)*(F2<260532)*(0.42*F2-8621.75)+(F2>260532)*(0.45*F2-16437.7)'
$ws.Range("L3").Formula = '=(C3>R3)*(C3<S3)*(974.58*(C3-R3)/10000+1400)*(C3-R3)/10000+(C3>S3)*(C3<T3)*((228.74*(C3-S3)/10000+2397)*(C3-S3)/10000+971)+(C3>T3)*(C3<U3)*(0.42*C3-8239)+(C3>U3)*(0.45*C3-15761)'
$ws.Range("M3").Formula = '=(D3>R3)*(D3<S3)*(974.58*(D3-R3)/10000+1400)*(D3-R3)/10000+(D3>S3)*(D3<T3)*((228.74*(D3-S3)/10000+2397)*(D3-S3)/10000+971)+(D3>T3)*(D3<U3)*(0.42*D3-8239)+(D3>U3)*(0.45*D3-15761)'
$ws.Range("N3").Formula = '=(E3>R3)*(E3<S3)*(974.58*(E3-R3)/10000+1400)*(E3-R3)/10000+(E3>S3)*(E3<T3)*((228.74*(E3-S3)/10000+2397)*(E3-S3)/10000+971)+(E3>T3)*(E3<U3)*(0.42*E3-8239)+(E3>U3)*(0.45*E3-15761)'
$ws.Range("O3").Formula = '=(F3>R3)*(F3<S3)*(974.58*(F3-R3)/10000+1400)*(F3-R3)/10000+(F3>S3)*(F3<U3)*((228.74*(F3-S3)/10000+2397)*(F3-T3)/10000+971)+(F3>U3)*(F3<U3)*(0.42*F3-8239)+(F3>U3)*(0.45*F3-15761)'
$ws.Range("L4").Formula = '=(C4>R4)*(C4<S4)*(912.17*(C4-R4)/10000+1400)*(C4-R4)/10000+(C4>S4)*(C4<T4)*((228.74*(C4-S4)/10000+2397)*(C4-S4)/10000+1038)+(C4>T4)*(C4<U4)*(0.42*C4-8172)+(C4>U4)*(0.45*C4-15694)'
$ws.Range("M4").Formula = '=(D4>R4)*(D4<S4)*(912.17*(D4-R4)/10000+1400)*(D4-R4)/10000+(D4>S4)*(D4<T4)*((228.74*(D4-S4)/10000+2397)*(D4-S4)/10000+1038)+(D4>T4)*(D4<U4)*(0.42*D4-8172)+(D4>U4)*(0.45*D4-15694)'
$ws.Range("N4").Formula = '=(E4>R4)*(E4<S4)*(912.17*(E4-R4)/10000+1400)*(E4-R4)/10000+(E4>S4)*(E4<T4)*((228.74*(E4-S4)/10000+2397)*(E4-S4)/10000+1038)+(E4>T4)*(E4<U4)*(0.42*E4-8172)+(E4>U4)*(0.45*E4-15694)'
$ws.Range("O4").Formula = '=(F4>R4)*(F4<S4)*(912.17*(F4-R4)/10000+1400)*(F4-R4)/10000+(F4>S4)*(F4<T4)*((228.74*(F4-S4)/10000+2397)*(F4-S4)/10000+1038)+(F4>T4)*(F4<U4)*(0.42*F4-8172)+(F4>U4)*(0.45*F4-15694)'
$ws.Range("L5").Formula = '=(C5>R5)*(C5<S5)*(883.74*(C5-R5)/10000+1500)*(C5-R5)/10000+(C5>S5)*(C5<T5)*((228.74*(C5-S5)/10000+2397)*(C5-S5)/10000+989)+(C5>T5)*(C5<U5)*(0.42*C5-7914)+(C5>U5)*(0.45*C5-15414)'
$ws.Range("M5").Formula = '=(D5>R5)*(D5<S5)*(883.74*(D5-R5)/10000+1500)*(D5-R5)/10000+(D5>S5)*(D5<T5)*((228.74*(D5-S5)/10000+2397)*(D5-S5)/10000+989)+(D5>T5)*(D5<U5)*(0.42*D5-7914)+(D5>U5)*(0.45*D5-15414)'
$ws.Range("N5").Formula = '=(E5>R5)*(E5<S5)*(883.74*(E5-R5)/10000+1500)*(E5-R5)/10000+(E5>S5)*(E5<T5)*((228.74*(E5-S5)/10000+2397)*(E5-S5)/10000+989)+(E5>T5)*(E5<U5)*(0.42*E5-7914)+(E5>U5)*(0.45*E5-15414)'
$ws.Range("O5").Formula = '=(F5>R5)*(F5<S5)*(883.74*(F5-R5)/10000+1500)*(F5-R5)/10000+(F5>S5)*(F5<T5)*((228.74*(F5-S5)/10000+2397)*(F5-S5)/10000+989)+(F5>T5)*(F5<U5)*(0.42*F5-7914)+(F5>U5)*(0.45*F5-15414)'
$ws.Range("L6").Formula = '=(C6>9000)*(C6<13996)*(997.8*(C6-9000)/10000+1400)*(C6-9000)/10000+(C6>13996)*(C6<54949)*((220.13*(C6-13996)/10000+2397)*(C6-13996)/10000+948.49)+(C6>54950)*(C6<260532)*(0.42*C6-8621.75)+(C6>260532)*(0.45*C6-16437.7)'
$ws.Range("M6").Formula = '=(D6>9000)*(D6<13996)*(997.8*(D6-9000)/10000+1400)*(D6-9000)/10000+(D6>13996)*(D6<54949)*((220.13*(D6-13996)/10000+2397)*(D6-13996)/10000+948.49)+(D6>54950)*(D6<260532)*(0.42*D6-8621.75)+(D6>260532)*(0.45*D6-16437.7)'
$ws.Range("N6").Formula = '=(E6>9000)*(E6<13996)*(997.8*(E6-9000)/10000+1400)*(E6-9000)/10000+(E6>13996)*(E6<54949)*((220.13*(E6-13996)/10000+2397)*(E6-13996)/10000+948.49)+(E6>54950)*(E6<260532)*(0.42*E6-8621.75)+(E6>260532)*(0.45*E6-16437.7)'
$ws.Range("O6").Formula = '=(F6>9000)*(F6<13996)*(997.8*(F6-9000)/10000+1400)*(F6-9000)/10000+(F6>13996)*(F6<54949)*((220.13*(F6-13996)/10000+2397)*(F6-13996)/10000+948.49)+(F6>54950)*(F6<260532)*(0.42*F6-8621.75)+(F6>260532)*(0.45*F6-16437.7)'

# --- Row 4 threshold corrections (R4/S4) ---
$ws.Range("R4").Value = 8005
$ws.Range("S4").Value = 13470

# --- View / selection update ---
# Scroll the window so column J is the leftmost visible column, then move
# the selection to N17 (mirrors <sheetView topLeftCell="J1"> / <selection
# activeCell="N17" sqref="N17"/> in the saved OOXML).
$win = $excel.ActiveWindow
$win.ScrollColumn = 10
$win.ScrollRow = 1
$ws.Range("N17").Select()
